$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# New reference-entry rows appended to "工作表1" (rows 55-60): apt-get, rfkill,
# lspci, ibus, tarball and patch notes - matching the commit "Add python and
# linux entries".
# ---------------------------------------------------------------------------

$c55 = @'
# Config of apt-get
There is config file of source provider address: /etc/apt/source.list, you can add your address if you are sure that the software you want to donwload is not available in public standard mirror. 
# Unmet Error
It is possible that adding new source address to this config file will cause **unmet** error when installing packages. This error message is like: Install xxx : Depends libxxx (version xxx) but version xxx is to be installed. To solve this problem, remove source address that caused error and run below statement to resume:
* apt-get clean
* apt-get autoclean
* apt-get update
'@

$c56 = @'
# What is rfkill
RF-kill is like an software level hardware switch. Say switch on/off the bluetooth/wireless service or others.
# Basic command
` sudo rfkill list all`
` sudo rfkill unblock all`
'@

$c57 = @'
# What is lspci
List PCI command list all hardware that is detected on PCI at hardware level (which means a hardware that appeared on lspci doesn’t mean it's driver is installed and can perform)
# Common usage
` lspci | grep -I network`  //show the network card
` lspci | grep -I ethernet`  //show the ethernet port
'@

$c58 = @'
# What is ibus
IBUS is a 輸入法 framework, it support chinese 倉頡, but make sure the system language installs chinese tradition at the first hand. Its on ubuntu but so far I only use the GUI of ibus. Ctrl-space to turn on an input method.
'@

$c59 = @'
 # Basic flow
* Download tarball: `wget {url} `
* Extract to /usr/local/src: `cd /usr/local/src; tar -zxvf ntp1.0.0.tar.gz`
* Config the make file: `cd ntp1.0.0/; ./configure --prefix=/usr/local/ntp`
* Make: `make clean; make; make check; make install`
* Link bin: can choose to build soft to /usr/local/bin or to add XXX_HOME and export XXX_HOME/bin to $PATH
* Link man: do something to MANPATH
# The Problem of arrangement
Consider a plan that set standard to install, bin putting and uninstall. If install without a clear domain, then its almost not possible to uninstall it later.
'@

$c60 = @'
# How to upgrade software wo patch
In the old tarball, use make to uninstall. Then download the new tarball, configure and make again. And still there can be a miss configure between 2 version (like when making the newer version, you forget to set some old setting)
# How patch help the upgrade
Use patch to update the source code, so that you dont need to configure the makefile again (keep the old config), but still need to make again or the software binaries will still not be updated/
'@

# Row 55 - apt-get / Source management
$ws.Cells.Item(55,1).Value = "apt-get"
$ws.Cells.Item(55,2).Value = "Source management"
$ws.Cells.Item(55,3).Value = $c55
$ws.Cells.Item(55,3).WrapText = $true

# Row 56 - rfkill / Basic
$ws.Cells.Item(56,1).Value = "rfkill"
$ws.Cells.Item(56,2).Value = "Basic"
$ws.Cells.Item(56,3).Value = $c56
$ws.Cells.Item(56,3).WrapText = $true

# Row 57 - lspci / List PCI Hardware
$ws.Cells.Item(57,1).Value = "lspci"
$ws.Cells.Item(57,2).Value = "List PCI Hardware"
$ws.Cells.Item(57,3).Value = $c57
$ws.Cells.Item(57,3).WrapText = $true

# Row 58 - ibus / A 輸入法 framework
$ws.Cells.Item(58,1).Value = "ibus"
$ws.Cells.Item(58,2).Value = "A 輸入法 framework"
$ws.Cells.Item(58,3).Value = $c58
$ws.Cells.Item(58,3).WrapText = $true

# Row 59 - tarball / Practice (with a red-highlighted "--prefix=..." run)
$ws.Cells.Item(59,1).Value = "tarball"
$ws.Cells.Item(59,2).Value = "Practice"
$ws.Cells.Item(59,3).Value = $c59
$ws.Cells.Item(59,3).WrapText = $true
$run59b = $ws.Cells.Item(59,3).Characters(176,24)
$run59b.Font.Color = 255
$run59b.Font.Name = "新細明體"
$run59c = $ws.Cells.Item(59,3).Characters(200,385)
$run59c.Font.Name = "新細明體"

# Row 60 - patch / Concept (with a red-highlighted "make" run)
$ws.Cells.Item(60,1).Value = "patch"
$ws.Cells.Item(60,2).Value = "Concept"
$ws.Cells.Item(60,3).Value = $c60
$ws.Cells.Item(60,3).WrapText = $true
$run60b = $ws.Cells.Item(60,3).Characters(428,4)
$run60b.Font.Color = 255
$run60b.Font.Name = "新細明體"
$run60c = $ws.Cells.Item(60,3).Characters(432,58)
$run60c.Font.Name = "新細明體"

# Reflect the post-edit selection/scroll position recorded in the workbook.
$ws.Range("C61").Select()
